# Updates Price (D) and Volume(1h) (E) columns with refreshed crypto quotes.
# For D-column values that look like plain numbers (e.g. "1.00", "19.37"),
# Excel's automatic type inference would otherwise coerce them into numeric
# cells and normalize their text (dropping trailing zeros, etc.). To keep
# them as literal text - matching the source data - we briefly force a text
# number format before assigning the value, then restore the cell style to
# "Normal" so no stray formatting is left behind on the cell itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.428.83"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "1.610.11"
$ws.Range("E3").Value = "  +1.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.60%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("D12").Value = "1.835.50"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").Value = "1.609.81"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "235.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.95%  "

$ws.Range("D18").Value = "26.441.28"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.91%  "

$ws.Range("D20").Value = "0.0₃0727"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.64%  "

$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.12%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("E31").Value = "  +0.13%  "

$ws.Range("D32").Value = "1.504.55"
$ws.Range("E32").Value = "  +5.93%  "

$ws.Range("E33").Value = "  +1.58%  "

$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("E36").Value = "  +1.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.935"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.22%  "

$ws.Range("D44").Value = "1.748.28"
$ws.Range("E44").Value = "  +1.15%  "

$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.42%  "

$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("E51").Value = "  +1.64%  "

